# Generate Report for Handback
#
# Refresh the handback-status report with the latest handoff/handback
# timestamps for the "abb6b744-1871-4545-a88f-6fac0f1d6f99.md" file:
#   - Overview sheet: "Latest HO Xliff Generate Date" column
#   - zh-cn / de-de sheets: "Correspond Handoff Datetime" and
#     "Correspond Handback DateTime" columns

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to abb6b744-1871-4545-a88f-6fac0f1d6f99.md
$wsOverview.Range("G3").Value = "2016-08-26 12:48:50"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 3 corresponds to abb6b744-1871-4545-a88f-6fac0f1d6f99.md
$wsZhCn.Range("H3").Value = "2016-08-26 12:48:46"
$wsZhCn.Range("K3").Value = "2016-08-26 12:49:07"

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Row 3 corresponds to abb6b744-1871-4545-a88f-6fac0f1d6f99.md
$wsDeDe.Range("H3").Value = "2016-08-26 12:48:50"
$wsDeDe.Range("K3").Value = "2016-08-26 12:49:16"
